$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
Write-Host ("Before H17=" + $ws.Cells.Item(17,8).Value)
$ws.Cells.Item(17,8).Value = 1903.5
Write-Host ("After H17=" + $ws.Cells.Item(17,8).Value)
